# v1.1 Updated the Testcases according to the review
# LH_TC_NOTIFICATION_Create_028
#
# Rename the first sheet from LH_TC_FEATURENAME to LH_TC_NOTIFICATION
# (this also re-targets the _xlnm._FilterDatabase defined name which
# refers to the sheet by name), and make that sheet the active tab
# (it was previously "VESRION HISTORY" that was active).

$wb = $excel.ActiveWorkbook

$wsNotification = $wb.Worksheets.Item("LH_TC_FEATURENAME")
$wsNotification.Name = "LH_TC_NOTIFICATION"

# Switch the active/selected tab from "VESRION HISTORY" to
# "LH_TC_NOTIFICATION" (first sheet).
$wsNotification.Activate()
